$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor-expressing-cell counts (1 -> 3) and the
# downstream expression/specificity statistics recomputed from them
# (columns F and L - the detection rates - are unchanged at 1).
# Column order matches: E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$cols = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)

$data = @(
    @(3,2.229205,6.687615,0.01491733042403661,0.01491733042403661,3,18.70849033333333,56.125471,0.8518406098705651,0.851840609870565,41.70506019351834,375.345541741665,0.01270718784605208,0.01270718784605208),
    @(3,2.229205,6.687615,0.01491733042403661,0.01491733042403661,3,1.675494666666667,5.026484,0.07628912719618224,0.07628912719618224,3.735021088406666,33.61518979566,0.001138030118146808,0.001138030118146808),
    @(3,2.229205,6.687615,0.01491733042403661,0.01491733042403661,3,1.578445666666666,4.735336999999999,0.07187026293325273,0.07187026293325273,3.518678972361666,31.668110751255,0.001072112459837722,0.001072112459837722),
    @(3,92.64273566666667,277.928207,0.6199440156137945,0.6199440156137944,3,18.70849033333333,56.125471,0.8518406098705651,0.851840609870565,1733.205724673389,15598.85152206049,0.5280934883460618,0.5280934883460617),
    @(3,92.64273566666667,277.928207,0.6199440156137945,0.6199440156137944,3,1.675494666666667,5.026484,0.07628912719618224,0.07628912719618224,155.2224095149098,1397.001685634188,0.04729498786167276,0.04729498786167275),
    @(3,92.64273566666667,277.928207,0.6199440156137945,0.6199440156137944,3,1.578445666666666,4.735336999999999,0.07187026293325273,0.07187026293325273,146.2315246611954,1316.083721950759,0.04455553940605995,0.04455553940605994),
    @(3,54.56532033333333,163.695961,0.365138653962169,0.365138653962169,3,18.70849033333333,56.125471,0.8518406098705651,0.851840609870565,1020.834767991403,9187.51291192263,0.3110399336784513,0.3110399336784512),
    @(3,54.56532033333333,163.695961,0.365138653962169,0.365138653962169,3,1.675494666666667,5.026484,0.07628912719618224,0.07628912719618224,91.42390320345821,822.815128831124,0.02785610921636268,0.02785610921636268),
    @(3,54.56532033333333,163.695961,0.365138653962169,0.365138653962169,3,1.578445666666666,4.735336999999999,0.07187026293325273,0.07187026293325273,86.12839343042853,775.1555408738568,0.02624261106735507,0.02624261106735507)
)

$startRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $cols[$j]).Value = $rowVals[$j]
    }
}